$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: "e" + "ditar" -> "editar" (keeps the first run's rPr) ---
$editar = $tr.Characters(15, 6)
$editar.Text = "editar"

# --- Paragraph 1: "s" + "ala" -> "sala" (keeps the second ("ala") run's rPr) ---
$sChar = $tr.Characters(22, 1)
$sChar.Text = ""
$alaRun = $tr.Characters(22, 3)
$alaRun.InsertBefore("s")

# --- Paragraph 2 (currently empty): add " Exámen y reportes" as 4 runs ---
$para2 = $tr.Paragraphs(2, 1)
$startPos = $para2.Start

$para2.InsertAfter(" ")
$para2.InsertAfter("Exámen")
$para2.InsertAfter(" y ")
$para2.InsertAfter("reportes")

$r1 = $tr.Characters($startPos, 1)
$r2 = $tr.Characters($startPos + 1, 6)
$r3 = $tr.Characters($startPos + 7, 3)
$r4 = $tr.Characters($startPos + 10, 8)

$r1.Font.Size = 32
$r2.Font.Size = 32
$r3.Font.Size = 32
$r4.Font.Size = 32
